$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2777502099621358
$ws.Cells.Item(2, 4).Value = 0.2504953273003991
$ws.Cells.Item(2, 5).Value = 0.3115690176358896
$ws.Cells.Item(2, 6).Value = 4.126984139852311
$ws.Cells.Item(2, 7).Value = 0.002377862666829956
$ws.Cells.Item(2, 9).Value = 0.5452418327394497
$ws.Cells.Item(2, 13).Value = 2.979882302059451
$ws.Cells.Item(3, 2).Value = 0.2443796282475432
$ws.Cells.Item(3, 4).Value = 0.2250659246027169
$ws.Cells.Item(3, 5).Value = 0.2722894332744659
$ws.Cells.Item(3, 6).Value = 3.756824166050166
$ws.Cells.Item(3, 7).Value = 0.00239504991606296
$ws.Cells.Item(3, 9).Value = 0.5233867455302814
$ws.Cells.Item(3, 13).Value = 2.622754573932099
$ws.Cells.Item(4, 2).Value = 0.223813822989797
$ws.Cells.Item(4, 4).Value = 0.2098915381554036
$ws.Cells.Item(4, 5).Value = 0.2481254865370346
$ws.Cells.Item(4, 6).Value = 3.535644766364783
$ws.Cells.Item(4, 7).Value = 0.002406064203868226
$ws.Cells.Item(4, 9).Value = 0.5104300635956704
$ws.Cells.Item(4, 13).Value = 2.4046568281438
$ws.Cells.Item(5, 2).Value = 0.2154142157494334
$ws.Cells.Item(5, 4).Value = 0.2038110996902844
$ws.Cells.Item(5, 5).Value = 0.2382656506253227
$ws.Cells.Item(5, 6).Value = 3.446956662512946
$ws.Cells.Item(5, 7).Value = 0.002410669711757161
$ws.Cells.Item(5, 9).Value = 0.5052664829403568
$ws.Cells.Item(5, 13).Value = 2.316044809461829
$ws.Cells.Item(6, 2).Value = 0.2140183361592278
$ws.Cells.Item(6, 4).Value = 0.2028074770275623
$ws.Cells.Item(6, 5).Value = 0.2366276232479407
$ws.Cells.Item(6, 6).Value = 3.432314745880888
$ws.Cells.Item(6, 7).Value = 0.00241144155740907
$ws.Cells.Item(6, 9).Value = 0.5044161074333218
$ws.Cells.Item(6, 13).Value = 2.301345920662698
$ws.Cells.Item(7, 2).Value = 0.223700618983429
$ws.Cells.Item(7, 4).Value = 0.2098091265554558
$ws.Cells.Item(7, 5).Value = 0.2479925668558849
$ws.Cells.Item(7, 6).Value = 3.534442955686416
$ws.Cells.Item(7, 7).Value = 0.002406125839903137
$ws.Cells.Item(7, 9).Value = 0.5103599543609718
$ws.Cells.Item(7, 13).Value = 2.403460745581043
$ws.Cells.Item(8, 2).Value = 0.266259973571465
$ws.Cells.Item(8, 4).Value = 0.2416316310054185
$ws.Cells.Item(8, 5).Value = 0.2980342055195564
$ws.Cells.Item(8, 6).Value = 3.998031708898964
$ws.Cells.Item(8, 7).Value = 0.002383693778048779
$ws.Cells.Item(8, 9).Value = 0.5376102371521085
$ws.Cells.Item(8, 13).Value = 2.856479783075486
$ws.Cells.Item(9, 2).Value = 0.3491090467171034
$ws.Cells.Item(9, 4).Value = 0.3078641449245936
$ws.Cells.Item(9, 5).Value = 0.3958651562032003
$ws.Cells.Item(9, 6).Value = 4.959786845367205
$ws.Cells.Item(9, 7).Value = 0.00234331229813355
$ws.Cells.Item(9, 9).Value = 0.5947221873233417
$ws.Cells.Item(9, 13).Value = 3.755782241541851
$ws.Cells.Item(10, 2).Value = 0.4096083368302743
$ws.Cells.Item(10, 4).Value = 0.3593597361082459
$ws.Cells.Item(10, 5).Value = 0.4676639709119002
$ws.Cells.Item(10, 6).Value = 5.704725183196729
$ws.Cells.Item(10, 7).Value = 0.002315768223600634
$ws.Cells.Item(10, 9).Value = 0.6389425668534727
$ws.Cells.Item(10, 13).Value = 4.425525274613022
$ws.Cells.Item(11, 2).Value = 0.4370524511855365
$ws.Cells.Item(11, 4).Value = 0.3835185285007867
$ws.Cells.Item(11, 5).Value = 0.5003363318206624
$ws.Cells.Item(11, 6).Value = 6.053395305479796
$ws.Cells.Item(11, 7).Value = 0.002303681874775688
$ws.Cells.Item(11, 9).Value = 0.6595572961037419
$ws.Cells.Item(11, 13).Value = 4.73272706042539
$ws.Cells.Item(12, 2).Value = 0.4474337798069996
$ws.Cells.Item(12, 4).Value = 0.3927824466056506
$ws.Cells.Item(12, 5).Value = 0.5127122600038945
$ws.Cells.Item(12, 6).Value = 6.186962634997428
$ws.Cells.Item(12, 7).Value = 0.002299167437556242
$ws.Cells.Item(12, 9).Value = 0.6674358625538304
$ws.Cells.Item(12, 13).Value = 4.849468714344511
$ws.Cells.Item(13, 2).Value = 0.4451984715187223
$ws.Cells.Item(13, 4).Value = 0.3907819922666818
$ws.Cells.Item(13, 5).Value = 0.5100466866984448
$ws.Cells.Item(13, 6).Value = 6.158126247283633
$ws.Cells.Item(13, 7).Value = 0.00230013694971491
$ws.Cells.Item(13, 9).Value = 0.665735851163987
$ws.Cells.Item(13, 13).Value = 4.824307252228465
$ws.Cells.Item(14, 2).Value = 0.4379067531076828
$ws.Cells.Item(14, 4).Value = 0.3842783019154297
$ws.Cells.Item(14, 5).Value = 0.5013544219129642
$ws.Cells.Item(14, 6).Value = 6.064352493581055
$ws.Cells.Item(14, 7).Value = 0.00230330922692873
$ws.Cells.Item(14, 9).Value = 0.6602040198724382
$ws.Cells.Item(14, 13).Value = 4.742322950202521
$ws.Cells.Item(15, 2).Value = 0.4334389110081531
$ws.Cells.Item(15, 4).Value = 0.3803099550984541
$ws.Cells.Item(15, 5).Value = 0.4960306940961487
$ws.Cells.Item(15, 6).Value = 6.007116908085038
$ws.Cells.Item(15, 7).Value = 0.002305260422346025
$ws.Cells.Item(15, 9).Value = 0.6568250351702858
$ws.Cells.Item(15, 13).Value = 4.692160165238107
$ws.Cells.Item(16, 2).Value = 0.4078132386920856
$ws.Cells.Item(16, 4).Value = 0.3577965043120059
$ws.Cells.Item(16, 5).Value = 0.4655291440028861
$ws.Cells.Item(16, 6).Value = 5.682146388557669
$ws.Cells.Item(16, 7).Value = 0.002316566903839923
$ws.Cells.Item(16, 9).Value = 0.6376054210648334
$ws.Cells.Item(16, 13).Value = 4.405503809860988
$ws.Cells.Item(17, 2).Value = 0.3920728687720043
$ws.Cells.Item(17, 4).Value = 0.3441800113659497
$ws.Cells.Item(17, 5).Value = 0.4468215939502755
$ws.Cells.Item(17, 6).Value = 5.485382812812304
$ws.Cells.Item(17, 7).Value = 0.002323615709730847
$ws.Cells.Item(17, 9).Value = 0.6259428553642152
$ws.Cells.Item(17, 13).Value = 4.230330414746675
$ws.Cells.Item(18, 2).Value = 0.383012146370362
$ws.Cells.Item(18, 4).Value = 0.3364165401638672
$ws.Cells.Item(18, 5).Value = 0.4360623801348567
$ws.Cells.Item(18, 6).Value = 5.373123558506393
$ws.Cells.Item(18, 7).Value = 0.002327711835398992
$ws.Cells.Item(18, 9).Value = 0.6192817544993261
$ws.Cells.Item(18, 13).Value = 4.129810685543902
$ws.Cells.Item(19, 2).Value = 0.3799430936331589
$ws.Cells.Item(19, 4).Value = 0.3337994078949009
$ws.Cells.Item(19, 5).Value = 0.4324195851671817
$ws.Cells.Item(19, 6).Value = 5.335267776282762
$ws.Cells.Item(19, 7).Value = 0.002329105940326385
$ws.Cells.Item(19, 9).Value = 0.6170344613623655
$ws.Cells.Item(19, 13).Value = 4.09581547886134
$ws.Cells.Item(20, 2).Value = 0.3937492117403565
$ws.Cells.Item(20, 4).Value = 0.3456223625037183
$ws.Cells.Item(20, 5).Value = 0.4488129396845153
$ws.Cells.Item(20, 6).Value = 5.506233156774101
$ws.Cells.Item(20, 7).Value = 0.002322861031921812
$ws.Cells.Item(20, 9).Value = 0.6271794982463348
$ws.Cells.Item(20, 13).Value = 4.248953231403021
$ws.Cells.Item(21, 2).Value = 0.4400488105831073
$ws.Cells.Item(21, 4).Value = 0.3861853754267202
$ws.Cells.Item(21, 5).Value = 0.5039074352543906
$ws.Cells.Item(21, 6).Value = 6.09185346576038
$ws.Cells.Item(21, 7).Value = 0.002302375770001725
$ws.Cells.Item(21, 9).Value = 0.6618268882764085
$ws.Cells.Item(21, 13).Value = 4.766392183978951
$ws.Cells.Item(22, 2).Value = 0.4702432314070109
$ws.Cells.Item(22, 4).Value = 0.4133740116275817
$ws.Cells.Item(22, 5).Value = 0.5399370367379248
$ws.Cells.Item(22, 6).Value = 6.483592040260305
$ws.Cells.Item(22, 7).Value = 0.002289350500272293
$ws.Cells.Item(22, 9).Value = 0.6848922070092414
$ws.Cells.Item(22, 13).Value = 5.10698805789562
$ws.Cells.Item(23, 2).Value = 0.4541338533372539
$ws.Cells.Item(23, 4).Value = 0.3987973739424433
$ws.Cells.Item(23, 5).Value = 0.5207045984905818
$ws.Cells.Item(23, 6).Value = 6.273646874802296
$ws.Cells.Item(23, 7).Value = 0.002296269588316413
$ws.Cells.Item(23, 9).Value = 0.6725430704706667
$ws.Cells.Item(23, 13).Value = 4.924967953460794
$ws.Cells.Item(24, 2).Value = 0.3929913725079928
$ws.Cells.Item(24, 4).Value = 0.3449700743720996
$ws.Cells.Item(24, 5).Value = 0.4479126648730585
$ws.Cells.Item(24, 6).Value = 5.496804038388632
$ws.Cells.Item(24, 7).Value = 0.002323202085453866
$ws.Cells.Item(24, 9).Value = 0.6266202754112413
$ws.Cells.Item(24, 13).Value = 4.240533268349168
$ws.Cells.Item(25, 2).Value = 0.3267614471486695
$ws.Cells.Item(25, 4).Value = 0.2894864509197248
$ws.Cells.Item(25, 5).Value = 0.3694210707617458
$ws.Cells.Item(25, 6).Value = 4.693344394412577
$ws.Cells.Item(25, 7).Value = 0.002353857934078113
$ws.Cells.Item(25, 9).Value = 0.5788768376679911
$ws.Cells.Item(25, 13).Value = 3.511091153429419
